$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row cell assignment (matches the appended rows 10-15 in the diff).
# Values that look numeric/date-like get a leading apostrophe so the engine
# stores them as text (shared string) instead of inferring number/date types,
# matching the target sharedStrings.xml / sheet1.xml content exactly.

# Row 10
$ws.Range("A10").Value = 'Vasya Pupkin'
$ws.Range("B10").Value = "'778821"
$ws.Range("C10").Value = "'10-05-2018"
$ws.Range("D10").Value = "'65874"
$ws.Range("E10").Value = 'Notebook'
$ws.Range("F10").Value = "'1"
$ws.Range("G10").Value = "'005478"
$ws.Range("H10").Value = 'Don''t power on'
$ws.Range("I10").Value = '+'

# Row 11
$ws.Range("A11").Value = 'Генрих'
$ws.Range("B11").Value = 'IV'
$ws.Range("C11").Value = "'03-06-1845"
$ws.Range("D11").Value = "'1"
$ws.Range("E11").Value = "'1"
$ws.Range("F11").Value = "'1"
$ws.Range("G11").Value = "'1"
$ws.Range("H11").Value = "'1"
$ws.Range("I11").Value = "'1"

# Row 12
$ws.Range("D12").Value = "'2"
$ws.Range("E12").Value = "'2"
$ws.Range("F12").Value = "'2"
$ws.Range("G12").Value = "'2"
$ws.Range("H12").Value = "'2"
$ws.Range("I12").Value = "'2"

# Row 13
$ws.Range("A13").Value = 'Адрон Солнцев'
$ws.Range("B13").Value = "'872645"
$ws.Range("C13").Value = "'07-09-2018"
$ws.Range("D13").Value = "'3"
$ws.Range("E13").Value = "'3"
$ws.Range("F13").Value = "'3"
$ws.Range("G13").Value = "'3"
$ws.Range("H13").Value = "'3"
$ws.Range("I13").Value = "'3"

# Row 14
$ws.Range("D14").Value = "'4"
$ws.Range("E14").Value = "'4"
$ws.Range("F14").Value = "'4"
$ws.Range("G14").Value = "'4"
$ws.Range("H14").Value = "'4"
$ws.Range("I14").Value = "'4"

# Row 15
$ws.Range("A15").Value = 'е'
$ws.Range("B15").Value = 'е'
$ws.Range("C15").Value = 'е'
$ws.Range("D15").Value = 'е'
$ws.Range("E15").Value = 'е'
$ws.Range("F15").Value = 'е'
$ws.Range("G15").Value = 'е'
$ws.Range("H15").Value = 'е'
$ws.Range("I15").Value = 'е'

# Reset the number format on the quote-prefixed cells back to the default
# "Normal" style so they do not carry an explicit style index (s="...")
# on the cell - matches styles.xml / sheet1.xml, where these cells have no
# style override, same as the rest of the sheet.
foreach ($addr in @("B10","C10","D10","F10","G10","C11","D11","E11","F11","G11","H11","I11","D12","E12","F12","G12","H12","I12","B13","C13","D13","E13","F13","G13","H13","I13","D14","E14","F14","G14","H14","I14")) {
    $ws.Range($addr).Style = "Normal"
}
